# weekly update, 2022-01-26, SARS-CoV-2
# Adds the newly released PDB entry 7ti9 (UBL1 domain of nsp3, form 2) to the
# "nsp3_Ubl1_" worksheet and registers it under the SARS-CoV-2 organism list
# on the "Organisms" worksheet.

$wb = $excel.ActiveWorkbook

# --- 1. nsp3_Ubl1_ sheet: append the new PDB entry as row 4 -----------------
$ws1 = $wb.Worksheets.Item("nsp3_Ubl1_")

$ws1.Cells.Item(4, 1).Value = "7ti9"
$ws1.Cells.Item(4, 2).Value = 2.73
$ws1.Cells.Item(4, 3).Value = "CRYSTAL STRUCTURE OF THE UBIQUITIN-LIKE DOMAIN 1 (UBL1) OF NSP3 FROM SARS-COV-2, FORM 2"
$ws1.Cells.Item(4, 4).Value = "X-RAY DIFFRACTION"

# The release date column looks like a date to Excel's autodetection, so it
# would silently be converted into a date serial number. Force it in as
# literal text (leading apostrophe), matching the source data which stores
# it as a plain string, then strip the "quote prefix" formatting that trick
# leaves behind so the cell keeps the sheet's default (unstyled) look.
$ws1.Cells.Item(4, 5).Value = "'2022-01-13"
$ws1.Cells.Item(3, 5).Copy()
$ws1.Cells.Item(4, 5).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Organisms sheet: insert 7ti9 into the SARS-CoV-2 PDB list ----------
$ws2 = $wb.Worksheets.Item("Organisms")

# Row 270 currently holds "6xg3" (the PDB id found at that offset within the
# "severe acute respiratory syndrome coronavirus2" block); insert a new row
# above it and push the remainder of the list down by one.
$ws2.Rows.Item(270).Insert()
$ws2.Cells.Item(270, 2).Value = "7ti9"

Write-Output "done"
